$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add two new columns P1, Q1 ---
# Copy the formatting of the existing last header cell (O1, style s="1":
# bold font, border, centered) onto the new header cells before setting
# their values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows (2-25): swap values in columns I/K/M/O and add P/Q ---
# For every data row the pattern is identical:
#   I: 1 -> 2
#   K: 2 -> 1
#   M: 1 -> 2
#   O: 2 -> 1
#   P: new column, value 2
#   Q: new column, value 2
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2
